$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '51.667.51'
Set-TextValue 'E2' '  -1.24%  '
Set-TextValue 'D3' '2.896.25'
Set-TextValue 'E3' '  -0.63%  '
Set-TextValue 'D5' '354.38'
Set-TextValue 'E5' '  +0.63%  '
Set-TextValue 'D6' '108.71'
Set-TextValue 'E6' '  -3.35%  '
Set-TextValue 'D7' '0.564'
Set-TextValue 'E7' '  +0.79%  '
Set-TextValue 'E8' '  +0.09%  '
Set-TextValue 'D9' '0.622'
Set-TextValue 'E9' '  -2.18%  '
Set-TextValue 'D10' '38.73'
Set-TextValue 'E10' '  -3.38%  '
Set-TextValue 'E11' '  +1.17%  '
Set-TextValue 'E12' '  +0.06%  '
Set-TextValue 'D13' '19.39'
Set-TextValue 'E13' '  -3.26%  '
Set-TextValue 'D14' '7.68'
Set-TextValue 'E14' '  -1.54%  '
Set-TextValue 'D15' '3.362.76'
Set-TextValue 'E15' '  -0.32%  '
Set-TextValue 'D16' '2.909.17'
Set-TextValue 'E16' '  -0.50%  '
Set-TextValue 'D17' '0.970'
Set-TextValue 'E17' '  -3.75%  '
Set-TextValue 'D18' '51.654.10'
Set-TextValue 'E18' '  -1.31%  '
Set-TextValue 'D19' '3.36'
Set-TextValue 'E19' '  +1.31%  '
Set-TextValue 'E20' '  -1.91%  '
Set-TextValue 'D21' '13.80'
Set-TextValue 'E21' '  -2.92%  '
Set-TextValue 'D22' '0.0₃0972'
Set-TextValue 'E22' '  -0.97%  '
Set-TextValue 'D23' '70.15'
Set-TextValue 'E23' '  -1.22%  '
Set-TextValue 'D24' '267.11'
Set-TextValue 'E24' '  -1.34%  '
Set-TextValue 'E25' '  -0.06%  '
Set-TextValue 'D26' '0.183'
Set-TextValue 'E26' '  +9.55%  '
Set-TextValue 'D27' '26.69'
Set-TextValue 'E27' '  -0.48%  '
Set-TextValue 'D28' '7.53'
Set-TextValue 'E28' '  +15.99%  '
Set-TextValue 'E29' '  +0.02%  '
Set-TextValue 'E30' '  +8.85%  '
Set-TextValue 'D31' '10.44'
Set-TextValue 'E31' '  -1.90%  '
Set-TextValue 'D32' '37.11'
Set-TextValue 'E32' '  -1.45%  '
Set-TextValue 'E33' '  -1.97%  '
Set-TextValue 'D34' '6.06'
Set-TextValue 'E34' '  -2.30%  '
Set-TextValue 'D35' '52.04'
Set-TextValue 'E35' '  -2.03%  '
Set-TextValue 'D36' '0.0438'
Set-TextValue 'E36' '  -2.79%  '
Set-TextValue 'E37' '  +0.07%  '
Set-TextValue 'E38' '  -4.51%  '
Set-TextValue 'D39' '18.11'
Set-TextValue 'E39' '  -3.95%  '
Set-TextValue 'D40' '1.98'
Set-TextValue 'E40' '  -4.02%  '
Set-TextValue 'D41' '2.68'
Set-TextValue 'E41' '  -7.28%  '
Set-TextValue 'D42' '0.118'
Set-TextValue 'E42' '  +0.89%  '
Set-TextValue 'D43' '22.80'
Set-TextValue 'E43' '  -3.51%  '
Set-TextValue 'D44' '118.88'
Set-TextValue 'E44' '  -1.85%  '
Set-TextValue 'D45' '2.17'
Set-TextValue 'E45' '  -0.78%  '
Set-TextValue 'D46' '2.47'
Set-TextValue 'E46' '  -6.09%  '
Set-TextValue 'D47' '3.41'
Set-TextValue 'E47' '  -4.13%  '
Set-TextValue 'D48' '2.119.27'
Set-TextValue 'E49' '  -6.11%  '
Set-TextValue 'D50' '0.0339'
Set-TextValue 'E50' '  +1.16%  '
Set-TextValue 'B51' 'FraxShare'
Set-TextValue 'C51' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D51' '9.02'
Set-TextValue 'E51' '  -0.64%  '
